# ProjectAdministration.xlsx edit script
# Commit: "Added some documentary comments, included identification rules and
#          fixed one rule that was wrong (project membership)"
#
# Net effect on worksheet data:
#  - The Dutch note "Niet meer in dienst." used for both A. Arends (C28) and
#    D. Diskstation (C31) is replaced by the English "has left the company".
#  - New documentary comments are added for C. Curly (C30) and O. Dysseus (C38).
#  - The D. Diskstation row (C31) additionally gets the plain/explicit cell
#    style that the rest of that column's data cells already use (fixing the
#    inconsistent "project membership" row formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Persons table: documentary comments / identification rules -----------

# A. Arends: translate the old "no longer employed" note to English.
$ws.Range("C28").Value = "has left the company"

# C. Curly: new comment - currently on maternity leave.
$ws.Range("C30").Value = "maternity leave"

# D. Diskstation: same translated note as A. Arends, and align its style
# with the rest of the column (this is the "project membership" rule fix).
$ws.Range("C31").Value = "has left the company"
$ws.Range("C31").Style = "Normal"

# O. Dysseus: new comment - on sabbatical.
$ws.Range("C38").Value = "is on sabbatical until 4th of July"

# --- View state: mirror where the author was last working -----------------

$ws.Range("C31").Select()
$win = $excel.ActiveWindow
try { $win.ScrollRow = 3 } catch {}
try { $win.ScrollColumn = 1 } catch {}
